$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing carbon-tax escalation formula in row 3: growth rate 1.1 -> 1.12 ---
$ws.Range("S3").Formula = "=R3*1.12^5"
$ws.Range("T3").Formula = "=S3*1.12^5"
$ws.Range("U3").Formula = "=T3*1.12^5"
$ws.Range("V3").Formula = "=U3*1.12^5"
$ws.Range("W3").Formula = "=V3*1.12^5"

# --- Row 4: CO2 / Process ---
$ws.Range("A4").Value = "CIMS.CAN"
$ws.Range("B4").Value = "Region"
$ws.Range("C4").Value = "CAN"
$ws.Range("G4").Value = "Tax"
$ws.Range("H4").Value = "CO2"

# --- Row 5: CH4 / Process ---
$ws.Range("A5").Value = "CIMS.CAN"
$ws.Range("B5").Value = "Region"
$ws.Range("C5").Value = "CAN"
$ws.Range("G5").Value = "Tax"
$ws.Range("H5").Value = "CH4"

# --- Sub_Context ("Process") for rows 4 and 5 ---
$ws.Range("I4").Value = "Process"
$ws.Range("I5").Value = "Process"

# --- Row 6: N2O / Process ---
$ws.Range("A6").Value = "CIMS.CAN"
$ws.Range("B6").Value = "Region"
$ws.Range("C6").Value = "CAN"
$ws.Range("G6").Value = "Tax"
$ws.Range("H6").Value = "N2O"
$ws.Range("I6").Value = "Process"

# --- Source / Unit + year values for row 4 ---
$ws.Range("K4").Value = "Govt of Canada"
$ws.Range("L4").Value = "$/tCO2e"
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Formula = "=S5"
$ws.Range("T4").Formula = "=T5"
$ws.Range("U4").Formula = "=U5"
$ws.Range("V4").Formula = "=V5"
$ws.Range("W4").Formula = "=W5"

# --- Source / Unit + year values for row 5 ---
$ws.Range("K5").Value = "Govt of Canada"
$ws.Range("L5").Value = "$/tCO2e"
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Formula = "=S3"
$ws.Range("T5").Formula = "=T3"
$ws.Range("U5").Formula = "=U3"
$ws.Range("V5").Formula = "=V3"
$ws.Range("W5").Formula = "=W3"

# --- Source / Unit + year values for row 6 ---
$ws.Range("K6").Value = "Govt of Canada"
$ws.Range("L6").Value = "$/tCO2e"
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Formula = "=S4"
$ws.Range("T6").Formula = "=T4"
$ws.Range("U6").Formula = "=U4"
$ws.Range("V6").Formula = "=V4"
$ws.Range("W6").Formula = "=W4"

# --- Update the selection to match the edited range (rows 3-6) ---
$ws.Range("A3:XFD6").Select()
